$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Price (column D) value (only where changed).
# These look numeric (e.g. "225.50", "10.02") but must stay literal text
# (matching the source sheet's inlineStr cells), so we force the cell to
# Text format before writing, then restore the default "Normal" style so
# no stray formatting is left behind.
$priceUpdates = [ordered]@{
    2  = "33.791.49"
    3  = "1.777.35"
    5  = "225.50"
    6  = "0.563"
    7  = "0.998"
    8  = "30.59"
    9  = "46.69"
    10 = "0.279"
    11 = "0.0666"
    12 = "0.0922"
    13 = "2.031.87"
    14 = "1.774.80"
    15 = "0.627"
    16 = "33.779.99"
    17 = "10.02"
    18 = "4.19"
    19 = "68.61"
    20 = "252.21"
    23 = "10.30"
    24 = "4.19"
    25 = "2.15"
    27 = "16.52"
    29 = "6.96"
    31 = "3.83"
    32 = "0.0514"
    35 = "1.85"
    36 = "1.483.85"
    37 = "1.07"
    39 = "83.33"
    42 = "2.69"
    43 = "0.886"
    45 = "0.0512"
    47 = "1.928.93"
    48 = "5.74"
    50 = "11.79"
    51 = "50.76"
}

# Map of row number -> new Volume(1h) (column E) value (all changed rows)
$volumeUpdates = [ordered]@{
    2  = "  +8.44%  "
    3  = "  +4.69%  "
    4  = "  -0.18%  "
    5  = "  +1.98%  "
    6  = "  +5.16%  "
    7  = "  -0.20%  "
    8  = "  +3.03%  "
    9  = "  +4.38%  "
    10 = "  +4.06%  "
    11 = "  +3.85%  "
    12 = "  +1.08%  "
    13 = "  +4.40%  "
    14 = "  +5.01%  "
    15 = "  +2.51%  "
    16 = "  +8.32%  "
    17 = "  -0.83%  "
    18 = "  +0.90%  "
    19 = "  +2.34%  "
    20 = "  +1.40%  "
    21 = "  +2.40%  "
    22 = "  -0.05%  "
    23 = "  +1.84%  "
    24 = "  -2.00%  "
    25 = "  -0.95%  "
    26 = "  +0.35%  "
    27 = "  +3.41%  "
    28 = "  +1.48%  "
    29 = "  +3.32%  "
    30 = "  -0.25%  "
    31 = "  +2.47%  "
    32 = "  +2.23%  "
    33 = "  +3.72%  "
    34 = "  +5.38%  "
    35 = "  +7.02%  "
    36 = "  -1.92%  "
    37 = "  +3.27%  "
    38 = "  +3.25%  "
    39 = "  +0.64%  "
    40 = "  +2.63%  "
    41 = "  +1.88%  "
    42 = "  +0.46%  "
    43 = "  +4.35%  "
    44 = "  +2.54%  "
    45 = "  +1.77%  "
    46 = "  +3.62%  "
    47 = "  +4.86%  "
    48 = "  +3.03%  "
    50 = "  +14.26%  "
    51 = "  -2.59%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
